$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (row 4): Video No. | Topic | concepts | explaination | command | example
$ws.Range("B4").Value = "Video No."
$ws.Range("C4").Value = "Topic"
$ws.Range("D4").Value = "concepts"
$ws.Range("E4").Value = "explaination"
$ws.Range("F4").Value = "command"
$ws.Range("G4").Value = "example"

# Clone the title cell's (B2) format onto the new header range first (so the
# range starts from an already-known style, i.e. bold/size 16) and only then
# shrink the font to 14 - this creates exactly one brand-new font entry
# (bold, size 14) instead of extra throw-away intermediate styles.
$ws.Range("B2").Copy()
$ws.Range("B4:G4").PasteSpecial(-4122)
$ws.Range("B4:G4").Font.Size = 14

# Row 4 should be a bit taller to fit the larger header font.
$ws.Rows.Item(4).RowHeight = 18

# Widen the columns used by the new table.
$ws.Range("B1:G1").ColumnWidth = 15.65

# Move the active selection to C11, matching the saved workbook state.
$ws.Range("C11").Select()
